# Auto-generated edit script: applies updated Leve profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR job sheets, as produced by
# the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 682.8
$ws.Range("I18").Value = 806.3333
$ws.Range("K18").Value = 806.3333
$ws.Range("M18").Value = -522.3333
# Row 32
$ws.Range("H32").Value = 797.6667
$ws.Range("J32").Value = 798.5
$ws.Range("L32").Value = 798.5
$ws.Range("N32").Value = -1450.5
# Row 33
$ws.Range("H33").Value = 137.7
$ws.Range("I33").Value = 125.22222
$ws.Range("K33").Value = 125.22222
$ws.Range("M33").Value = 103.77778
# Row 39
$ws.Range("H39").Value = 179.41176
$ws.Range("I39").Value = 193
$ws.Range("J39").Value = 116
$ws.Range("K39").Value = 579
$ws.Range("L39").Value = 348
$ws.Range("M39").Value = -283
$ws.Range("N39").Value = -940
# Row 40
$ws.Range("H40").Value = 6342
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10349
# Row 70
$ws.Range("H70").Value = 4838.2
$ws.Range("I70").Value = 2128.3333
$ws.Range("K70").Value = 6384.999899999999
$ws.Range("M70").Value = -6114.999899999999
# Row 73
$ws.Range("H73").Value = 4838.2
$ws.Range("I73").Value = 2128.3333
$ws.Range("K73").Value = 6384.999899999999
$ws.Range("M73").Value = -5448.999899999999
# Row 113
$ws.Range("H113").Value = 1692.1666
$ws.Range("J113").Value = 991.25
$ws.Range("L113").Value = 991.25
$ws.Range("N113").Value = -7499.25
# Row 125
$ws.Range("H125").Value = 3416
$ws.Range("J125").Value = 3433.3333
$ws.Range("L125").Value = 30899.9997
$ws.Range("N125").Value = -35819.9997
# Row 132
$ws.Range("H132").Value = 11077.954
$ws.Range("I132").Value = 13732.25
$ws.Range("K132").Value = 41196.75
$ws.Range("M132").Value = -38666.75

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 19999.5
$ws.Range("J23").Value = 19999.5
$ws.Range("L23").Value = 19999.5
$ws.Range("N23").Value = -20517.5
# Row 45
$ws.Range("H45").Value = 2448.1428
$ws.Range("J45").Value = 3507.7144
$ws.Range("L45").Value = 3507.7144
$ws.Range("N45").Value = -4261.7144

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 134
$ws.Range("H134").Value = 2583
$ws.Range("I134").Value = 1777.3334
$ws.Range("K134").Value = 5332.0002
$ws.Range("M134").Value = -2797.0002

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5541.5713
$ws.Range("I58").Value = 4959.2
$ws.Range("K58").Value = 4959.2
$ws.Range("M58").Value = -4756.2
# Row 136
$ws.Range("H136").Value = 5541.5713
$ws.Range("I136").Value = 4959.2
$ws.Range("K136").Value = 14877.6
$ws.Range("M136").Value = -12327.6

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 258.16666
$ws.Range("J23").Value = 258.16666
$ws.Range("L23").Value = 774.4999799999999
$ws.Range("N23").Value = -1244.49998
# Row 34
$ws.Range("H34").Value = 929.2
$ws.Range("J34").Value = 1188.5454
$ws.Range("L34").Value = 3565.6362
$ws.Range("N34").Value = -3733.6362
# Row 113
$ws.Range("H113").Value = 1509.25
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 1712.3334
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 5137.0002
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -9477.0002

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 11978280
$ws.Range("I11").Value = 11347850
$ws.Range("K11").Value = 11347850
$ws.Range("M11").Value = -11347711
# Row 29
$ws.Range("H29").Value = 18319.8
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 22774.75
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 22774.75
$ws.Range("M29").Value = -210
$ws.Range("N29").Value = -23354.75
# Row 97
$ws.Range("H97").Value = 716.2727
$ws.Range("J97").Value = 874.8
$ws.Range("L97").Value = 874.8
$ws.Range("N97").Value = -1866.8

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7973.6
$ws.Range("I7").Value = 7973.6
$ws.Range("K7").Value = 7973.6
$ws.Range("M7").Value = -7861.6
# Row 14
$ws.Range("H14").Value = 10900
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10900
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = -11244
$ws.Range("L14").Value = 10900
$ws.Range("M14").ClearContents()
# Row 46
$ws.Range("H46").Value = 5744.95
$ws.Range("I46").Value = 4233.1665
$ws.Range("J46").Value = 6392.857
$ws.Range("K46").Value = 4233.1665
$ws.Range("L46").Value = 6392.857
$ws.Range("M46").Value = -4045.1665
$ws.Range("N46").Value = -6768.857
# Row 55
$ws.Range("H55").Value = 1389.3
$ws.Range("I55").Value = 1267.875
$ws.Range("J55").Value = 1875
$ws.Range("K55").Value = 1267.875
$ws.Range("L55").Value = 1875
$ws.Range("M55").Value = -1094.875
$ws.Range("N55").Value = -2221
# Row 74
$ws.Range("H74").Value = 32722
$ws.Range("I74").Value = 15444
$ws.Range("K74").Value = 15444
$ws.Range("M74").Value = -14446
# Row 77
$ws.Range("H77").Value = 32722
$ws.Range("I77").Value = 15444
$ws.Range("K77").Value = 46332
$ws.Range("M77").Value = -41340
# Row 126
$ws.Range("H126").Value = 7973.6
$ws.Range("I126").Value = 7973.6
$ws.Range("K126").Value = 23920.8
$ws.Range("M126").Value = -21450.8

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 2500007.5
$ws.Range("I21").Value = 2500007.5
$ws.Range("K21").Value = 2500007.5
$ws.Range("M21").Value = -2499772.5
# Row 35
$ws.Range("H35").Value = 2500007.5
$ws.Range("I35").Value = 2500007.5
$ws.Range("K35").Value = 2500007.5
$ws.Range("M35").Value = -2499717.5
# Row 41
$ws.Range("H41").Value = 34812.25
$ws.Range("J41").Value = 34544.668
$ws.Range("L41").Value = 34544.668
$ws.Range("N41").Value = -35324.668
# Row 51
$ws.Range("H51").Value = 44999.5
$ws.Range("I51").Value = 44999.5
$ws.Range("K51").Value = 44999.5
$ws.Range("M51").Value = -44489.5
# Row 58
$ws.Range("H58").Value = 18000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 69
$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 2500
$ws.Range("L69").Value = 2500
$ws.Range("N69").Value = -3998
# Row 72
$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 2500
$ws.Range("L72").Value = 7500
$ws.Range("N72").Value = -14988
# Row 136
$ws.Range("H136").Value = 3183.6
$ws.Range("I136").Value = 1980.75
$ws.Range("K136").Value = 5942.25
$ws.Range("M136").Value = -3392.25

